$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1906.8036
$ws.Range("J17").Value = 1906.8036
$ws.Range("L17").Value = 5720.4108
$ws.Range("N17").Value = -6056.4108

$ws.Range("H127").Value = 2234.25
$ws.Range("I127").Value = 1606.6666
$ws.Range("J127").Value = 4117
$ws.Range("K127").Value = 4819.9998
$ws.Range("L127").Value = 12351
$ws.Range("M127").Value = 140.0002000000004
$ws.Range("N127").Value = -22271

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 68685
$ws.Range("J131").Value = 68685
$ws.Range("L131").Value = 68685
$ws.Range("N131").Value = -78765

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").ClearContents()
$ws.Range("N6").Value = 0

$ws.Range("H22").Value = 363.875
$ws.Range("I22").Value = 374.85715
$ws.Range("K22").Value = 374.85715
$ws.Range("M22").Value = -201.85715

$ws.Range("H86").Value = 2688.3333
$ws.Range("I86").Value = 2928
$ws.Range("J86").Value = 1849.5
$ws.Range("K86").Value = 2928
$ws.Range("L86").Value = 1849.5
$ws.Range("M86").Value = -1805
$ws.Range("N86").Value = -4095.5

$ws.Range("H89").Value = 2688.3333
$ws.Range("I89").Value = 2928
$ws.Range("J89").Value = 1849.5
$ws.Range("K89").Value = 14640
$ws.Range("L89").Value = 9247.5
$ws.Range("M89").Value = -9024
$ws.Range("N89").Value = -20479.5

$ws.Range("H94").Value = 395.7143
$ws.Range("I94").Value = 395.7143
$ws.Range("K94").Value = 395.7143
$ws.Range("M94").Value = 55.28570000000002

$ws.Range("H95").Value = 17000
$ws.Range("J95").Value = 17000
$ws.Range("L95").Value = 17000
$ws.Range("N95").Value = -22492

$ws.Range("H97").Value = 15555
$ws.Range("I97").Value = 15555
$ws.Range("K97").Value = 15555
$ws.Range("M97").Value = -14564

$ws.Range("H99").Value = 2502.5

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3067.3333
$ws.Range("I62").Value = 2601.5
$ws.Range("J62").Value = 3999
$ws.Range("K62").Value = 2601.5
$ws.Range("L62").Value = 3999
$ws.Range("M62").Value = -1977.5
$ws.Range("N62").Value = -5247

$ws.Range("H65").Value = 3067.3333
$ws.Range("I65").Value = 2601.5
$ws.Range("J65").Value = 3999
$ws.Range("K65").Value = 13007.5
$ws.Range("L65").Value = 19995
$ws.Range("M65").Value = -9887.5
$ws.Range("N65").Value = -26235

$ws.Range("H74").Value = 49000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 49000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H86").Value = 7750
$ws.Range("J86").Value = 7750
$ws.Range("L86").Value = 7750
$ws.Range("N86").Value = -9996

$ws.Range("H89").Value = 7750
$ws.Range("J89").Value = 7750
$ws.Range("L89").Value = 38750
$ws.Range("N89").Value = -49982

$ws.Range("H96").Value = 23055.25
$ws.Range("J96").Value = 23055.25
$ws.Range("L96").Value = 23055.25
$ws.Range("N96").Value = -28547.25

$ws.Range("H122").Value = 3728.4546
$ws.Range("I122").Value = 3999.8572
$ws.Range("K122").Value = 11999.5716
$ws.Range("M122").Value = -9549.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1257.8
$ws.Range("I5").Value = 1257.8
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3773.4
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -3661.4

$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("K16").Value = 300
$ws.Range("M16").Value = -127

$ws.Range("H34").Value = 52332.43
$ws.Range("J34").Value = 60940.332
$ws.Range("L34").Value = 182820.996
$ws.Range("N34").Value = -182988.996

$ws.Range("H39").Value = 9000
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H55").Value = 13532.333
$ws.Range("J55").Value = 19998.5
$ws.Range("L55").Value = 59995.5
$ws.Range("N55").Value = -60349.5

$ws.Range("H75").Value = 10013
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 10013
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H92").Value = 225.42857
$ws.Range("I92").Value = 250.6
$ws.Range("J92").Value = 162.5
$ws.Range("K92").Value = 751.8
$ws.Range("L92").Value = 487.5
$ws.Range("M92").Value = 496.2
$ws.Range("N92").Value = -2983.5

$ws.Range("H135").Value = 1257.8
$ws.Range("I135").Value = 1257.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 11320.2
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -8785.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2133.3333
$ws.Range("I102").Value = 2133.3333
$ws.Range("K102").Value = 2133.3333
$ws.Range("M102").Value = -511.3332999999998

$ws.Range("H128").Value = 37450
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 663.3
$ws.Range("I46").Value = 692
$ws.Range("J46").Value = 596.3333
$ws.Range("K46").Value = 692
$ws.Range("L46").Value = 596.3333
$ws.Range("M46").Value = -504
$ws.Range("N46").Value = -972.3333

$ws.Range("H98").Value = 63999
$ws.Range("J98").Value = 63999
$ws.Range("L98").Value = 63999
$ws.Range("N98").Value = -69989

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6665.727
$ws.Range("I122").Value = 7732.5713
$ws.Range("J122").Value = 4798.75
$ws.Range("K122").Value = 23197.7139
$ws.Range("L122").Value = 14396.25
$ws.Range("M122").Value = -20747.7139
$ws.Range("N122").Value = -19296.25

$ws.Range("H130").Value = 40997.5
$ws.Range("J130").Value = 40997.5
$ws.Range("L130").Value = 40997.5
$ws.Range("N130").Value = -51037.5
